$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Fix the placeholder row 28 entry and add a new row 27 for two new books.

# Row 27: Rich Dad Poor Dad / Robert Kiyosaki / Работа / 2. Читаю
$ws.Range("A27").Value = "Rich Dad Poor Dad"
$ws.Range("B27").Value = "Robert Kiyosaki"
$ws.Range("C27").Value = "Работа"
$ws.Range("F27").Value = "2. Читаю"

# Row 28: The Richest Man in Babylon / George Samuel Clason / Работа / 1. В очереди
$ws.Range("A28").Value = "The Richest Man in Babylon"
$ws.Range("B28").Value = "George Samuel Clason"
$ws.Range("C28").Value = "Работа"
$ws.Range("F28").Value = "1. В очереди"

$ws.Range("G20").Select()
